$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 586, shifting existing rows 586-627 down to 587-628
$ws.Rows.Item(586).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds dates as plain text (matching the rest of the sheet), so
# force text formatting before assigning; otherwise Excel auto-converts a
# recognizable date string like "2026/01/07" into a real date serial value.
$ws.Cells.Item(586, 1).NumberFormat = "@"
$ws.Cells.Item(586, 1).Value = "2026/01/07"
$ws.Cells.Item(586, 1).ClearFormats()

$ws.Cells.Item(586, 2).Value = "水"
$ws.Cells.Item(586, 3).Value = 17
$ws.Cells.Item(586, 4).Value = 201
